$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-09-09"
$ws.Range("I1").Value = "2022 (through 09-09)"
$ws.Range("I10").Value = 39
$ws.Range("I14").Value = 1176
